# Append a new trade row (row 8) to the sheet, mirroring the existing
# data rows (columns A:H).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 1).Value = 9565.85
$ws.Cells.Item(8, 2).Value = 9666.3799999999992
$ws.Cells.Item(8, 3).Value = 109.08
$ws.Cells.Item(8, 4).Value = 107.95
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = -1.04
$ws.Cells.Item(8, 7).Value = 42612.672974537039
$ws.Cells.Item(8, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(8, 8).Value = $false
